$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.588.88"
$ws.Range("E2").Value = "  -4.40%  "
$ws.Range("D3").Value = "3.146.20"
$ws.Range("E3").Value = "  -4.64%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "524.36"
$ws.Range("E5").Value = "  -6.17%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "133.07"
$ws.Range("E6").Value = "  -6.94%  "
$ws.Range("E7").Value = "  -0.12%  "
$ws.Range("D8").Value = "3.148.33"
$ws.Range("E8").Value = "  -4.68%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.452"
$ws.Range("E9").Value = "  -6.02%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.27"
$ws.Range("E10").Value = "  -7.27%  "
$ws.Range("E11").Value = "  -6.84%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.389"
$ws.Range("E12").Value = "  -4.09%  "
$ws.Range("D13").Value = "3.693.90"
$ws.Range("E13").Value = "  -4.53%  "
$ws.Range("E14").Value = "  -1.72%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "25.70"
$ws.Range("E15").Value = "  -5.50%  "
$ws.Range("D16").Value = "3.148.99"
$ws.Range("E16").Value = "  -4.64%  "
$ws.Range("D17").Value = "57.602.36"
$ws.Range("E17").Value = "  -4.44%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.0000152"
$ws.Range("E18").Value = "  -8.00%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.80"
$ws.Range("E19").Value = "  -5.20%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.03"
$ws.Range("E20").Value = "  -8.89%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.03"
$ws.Range("E21").Value = "  -6.05%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "344.94"
$ws.Range("E22").Value = "  -7.67%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("E23").Value = "  -0.25%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "69.51"
$ws.Range("E24").Value = "  -5.93%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.510"
$ws.Range("E25").Value = "  -6.57%  "
$ws.Range("D26").Value = "3.283.66"
$ws.Range("E26").Value = "  -4.96%  "
$ws.Range("D27").Value = "0.0₃0951"
$ws.Range("E27").Value = "  -8.81%  "
$ws.Range("E28").Value = "  -4.97%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.996"
$ws.Range("E29").Value = "  -0.11%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.83"
$ws.Range("E30").Value = "  -4.96%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.997"
$ws.Range("E31").Value = "  -0.26%  "
$ws.Range("E32").Value = "  -8.00%  "
$ws.Range("E33").Value = "  -8.91%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "21.67"
$ws.Range("E34").Value = "  -3.73%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.20"
$ws.Range("E35").Value = "  -5.54%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.90"
$ws.Range("E36").Value = "  -5.36%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "159.08"
$ws.Range("E37").Value = "  -4.65%  "
$ws.Range("E38").Value = "  -6.99%  "
$ws.Range("E39").Value = "  -6.93%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "25.71"
$ws.Range("E40").Value = "  -5.26%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0695"
$ws.Range("E41").Value = "  -5.10%  "
$ws.Range("D42").Value = "3.178.58"
$ws.Range("E42").Value = "  -4.76%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "40.60"
$ws.Range("E43").Value = "  -3.20%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.695"
$ws.Range("E44").Value = "  -7.22%  "
$ws.Range("E45").Value = "  -4.32%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.94"
$ws.Range("E46").Value = "  -5.82%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.00"
$ws.Range("E47").Value = "  -0.09%  "
$ws.Range("E48").Value = "  -7.73%  "
$ws.Range("D49").Value = "2.262.32"
$ws.Range("E49").Value = "  -4.04%  "
$ws.Range("E50").Value = "  -5.44%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "20.48"
$ws.Range("E51").Value = "  -4.17%  "
